# Update Equipment Load Factors values per RMI data refresh.

$wb = $excel.ActiveWorkbook

# --- ELF-bldg-winter ---
$ws = $wb.Worksheets.Item("ELF-bldg-winter")
$ws.Range("B2").Value = 1.36302   # heating, urban residential
$ws.Range("D2").Value = 1.16033   # heating, commercial
$ws.Range("B5").Value = 1.91535   # lighting, urban residential
$ws.Range("D5").Value = 1.45616   # lighting, commercial
$ws.Range("D7").Value = 1.45616   # other, commercial

# --- ELF-bldg-summer ---
$ws = $wb.Worksheets.Item("ELF-bldg-summer")
$ws.Range("B3").Value = 7.66676   # cooling & ventilation, urban residential
$ws.Range("D3").Value = 6.54006   # cooling & ventilation, commercial
$ws.Range("B5").Value = 1.85042   # lighting, urban residential
$ws.Range("D5").Value = 2.00709   # lighting, commercial
$ws.Range("D7").Value = 2.00709   # other, commercial

# --- ELF-vehicles ---
$ws = $wb.Worksheets.Item("ELF-vehicles")
$ws.Range("B4").Value = 1.16038   # aircraft, summer
$ws.Range("C4").Value = 1.22331   # aircraft, winter
$ws.Range("B5").Value = 1.16038   # rail, summer
$ws.Range("C5").Value = 1.22331   # rail, winter
$ws.Range("B6").Value = 1.16038   # ships, summer
$ws.Range("C6").Value = 1.22331   # ships, winter
$ws.Range("B7").Value = 1.16038   # motorbikes, summer
$ws.Range("C7").Value = 1.22331   # motorbikes, winter
